# Add eight new metadata columns to the "Tabelle5" table on the
# "Booklet_FK Lagerlogistik" sheet, and populate the first of them
# ("AssessmentType") with a default value of 0 for every existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

$newColumnNames = @(
    "AssessmentType",
    "Description",
    "Disclaimer",
    "Duration",
    "EscoOccupationId",
    "EscoSkills",
    "Publisher",
    "Title"
)

# Grab the formatting of the current last header cell so the new headers
# match the look of the existing ones (bold header row style).
$lastHeaderCell = $lo.HeaderRowRange.Cells.Item(1, $lo.HeaderRowRange.Columns.Count)

$firstNewColumn = $null

foreach ($name in $newColumnNames) {
    $col = $lo.ListColumns.Add()
    $headerCell = $col.Range.Item(1)
    $headerCell.Value = $name

    $lastHeaderCell.Copy() | Out-Null
    $headerCell.PasteSpecial(-4122) | Out-Null

    if ($null -eq $firstNewColumn) {
        $firstNewColumn = $col
    }
}

$excel.CutCopyMode = 0

# The first new column ("AssessmentType") gets a default value of 0 for
# every existing data row.
$firstNewColumn.DataBodyRange.Value = 0

# Reflect the interactive selection/scroll position left behind by the edit.
$ws.Range("BJ2").Select() | Out-Null
